# Applies the "Ceri work from teams" edit to minutesweek3.2.docx
$d = $word.ActiveDocument

# -------------------------------------------------------------------
# 1) Date paragraph (paragraph 3): "Date:     14/02/2023 4.30 -5.00pm"
#    -> "Date and time:     4.30 -5.00pm, 14/02/2023"
# -------------------------------------------------------------------
$p = $d.Paragraphs(3)
$r = $d.Range($p.Range.Start, $p.Range.End)
$r.Text = "Date and time:     4.30 -5.00pm, 14/02/2023"

Write-Output "Step1 done: $($d.Paragraphs(3).Range.Text)"

# -------------------------------------------------------------------
# 2) Paragraph 7 ("Tasks of this meeting:") becomes the (edited) body
#    text that used to live in paragraph 9.
# -------------------------------------------------------------------
$p7 = $d.Paragraphs(7)
$r7 = $d.Range($p7.Range.Start, $p7.Range.End)
$r7.Text = "On this meeting, taken straight after the interview, we discussed how the interview went. " + `
           "We decided to make a web-based software. " + `
           "We had a read through of the documents added by lecturer and decided what " + `
           "next steps we take. The documents we need to prepare, which we are aware " + `
           "of now, are Requirements Specification document, UML and Use Case diagram. "

Write-Output "Step2 done: $($d.Paragraphs(7).Range.Text)"

# -------------------------------------------------------------------
# 3) Paragraph 9 (formerly "On this meeting...") becomes the new
#    "Next meeting: 2 - 3pm on Monday, 20th of February" line, with
#    "th" as a superscript run.
# -------------------------------------------------------------------
$p9 = $d.Paragraphs(9)
$p9Start = $p9.Range.Start
$r9 = $d.Range($p9Start, $p9.Range.End)
$enDash = [char]0x2013
$r9.Text = "Next meeting: 2 $enDash 3pm on Monday, 20th of February "

$full9 = $d.Paragraphs(9).Range.Text
$thIdx = $full9.IndexOf("20th") + 2
$thRange = $d.Range($p9Start + $thIdx, $p9Start + $thIdx + 2)
$thRange.Font.Superscript = $true

Write-Output "Step3 done: $($d.Paragraphs(9).Range.Text)"

# -------------------------------------------------------------------
# 4) Collapse the two blank paragraphs (10 & 11) that used to sit
#    between the old "Requirements..." paragraph and "Action for the
#    next meeting:" down to a single blank paragraph.
# -------------------------------------------------------------------
$d.Paragraphs(10).Range.Delete()

Write-Output "Step4 done. Count=$($d.Paragraphs.Count)"

# -------------------------------------------------------------------
# 5) Remove the blank paragraph between "Action for the next
#    meeting:" and the "Secretary sends..." bullet.
# -------------------------------------------------------------------
$d.Paragraphs(12).Range.Delete()

Write-Output "Step5 done. Count=$($d.Paragraphs.Count)"

# -------------------------------------------------------------------
# 6) "...Sam is taking the lead; on next meeting..." -> "...lead, on
#    next meeting..." (semicolon -> comma, runs merged).
# -------------------------------------------------------------------
$p13 = $d.Paragraphs(13)
$r13 = $d.Range($p13.Range.Start, $p13.Range.End)
$r13.Text = "We start working on the requirements specification document. Sam is taking the lead, on next meeting we will analyse prepared text and take further steps. "

Write-Output "Step6 done: $($d.Paragraphs(13).Range.Text)"

# -------------------------------------------------------------------
# 7) Collapse the two blank paragraphs (15 & 16) after the bullet
#    list down to a single blank paragraph.
# -------------------------------------------------------------------
$d.Paragraphs(16).Range.Delete()

Write-Output "Step7 done. Count=$($d.Paragraphs.Count)"

# -------------------------------------------------------------------
# 8) The old "Next meeting: Tuesday, 21/02/2023 3pm-4pm" paragraph
#    (now paragraph 16) becomes two space-only runs, and the two
#    trailing blank paragraphs after it are removed.
# -------------------------------------------------------------------
$d.Paragraphs(18).Range.Delete()
$d.Paragraphs(17).Range.Delete()

$pLast = $d.Paragraphs(16)
$pLastStart = $pLast.Range.Start
$rLast = $d.Range($pLastStart, $pLast.Range.End)
$rLast.Text = " "
$afterFirst = $d.Range($pLastStart + 1, $pLastStart + 1)
$afterFirst.InsertAfter(" ")

Write-Output "Step8 done. Count=$($d.Paragraphs.Count)"
for ($i=1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Output "  $i : [$($d.Paragraphs($i).Range.Text)]"
}
